# Expanded the error-handling of the excel importer
# Insert a new blank row at the top of the data sheet (Blad1 / sheet1),
# shifting all existing rows (header + 171 data rows) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 1; all rows shift down by one.
$ws.Rows.Item(1).Insert()

# Update the visible selection to match the post-edit view (top-left back
# to A1, active cell F5).
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("F5").Select()
